# Commit: "update file with jgit"
# The "Rules" sheet cell E8 previously held the shared string "Good Morning".
# It is updated to "GIT UPDATE", and that cell becomes the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select() | Out-Null
